$d = $word.ActiveDocument

# Locate the existing "_GoBack" bookmark and remove it from its current
# location (it needs to move to the very end of the newly-added content).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Work from the very end of the document (the last paragraph is
# "För att veta vad jag ska göra ... exjobb)").
$tailRange = $d.Paragraphs.Last.Range

# 1) New empty "No Spacing" paragraph.
$tailRange.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("No Spacing")

# 2) "Dag 6" heading.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Heading 3")
$p.Range.InsertAfter("Dag 6")

# 3) Paragraph about Åke's feedback.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.InsertAfter("Idag har jag fått påpekat av Åke att det är väldigt dålig approach om jag skulle göra klart mitt formulär för att sedan hitta de artiklar som är relevanta. Egentligen borde jag hittat dessa artiklar innan jag börjat för att utveckla efter vad artiklarna säger. Så det jag får göra nu är att utveckla under tiden jag hittar artiklar. ")

# 4) Paragraph about the article and Jackson.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.InsertAfter("Jag hittade en artikel som i princip täcker upp alla delar jag borde undersöka gällande designen och funktionen av formuläret. Jag pratade med Jackson idag och han hade lite mycket att göra så jag skulle återkomma nästa vecka för att få hjälp att sätta upp en server och få tillgång till denna. ")

# 5) Paragraph about the database relations.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.InsertAfter("Gällande databasen har jag börjat fundera på relationerna mellan tabellerna. Jag har skapat en tabell som håller ihop alla andra tabeller för att få en grupp kopplat till en ansökning som man kan hämta i översiktten. Men ska t.ex en person kunna göra flera ansökningar med samma uppgifter?")

# 6) Paragraph about the articles to look for; this is where the
#    "_GoBack" bookmark ends up, right at the end of the inserted text.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.InsertAfter("Artiklar som jag kommer leta efter är om design för formuläret, hur jag borde genomföra mitt test och ta vara på datat samt hur ställer sig formuläret mot dem frågor som har varit ofullständiga")

# Placing a bookmark collapsed exactly at "end of paragraph text" (one
# position before the paragraph mark) trips an edge case in this COM
# runtime, so append a throwaway sentinel character first, anchor the
# bookmark just before it, then remove the sentinel again.
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("X")
$p = $d.Paragraphs.Last
$sentinelPos = $p.Range.End - 1
$bmPos = $sentinelPos - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($bmPos, $sentinelPos).Delete()
